# Updated symbol list on Tue Dec 27 08:32:37 UTC 2022 with GitHub Actions
# Applies the refreshed cryptocurrency price/volume snapshot to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text cells (coin name, link, volume label) - safe to set directly.
$textValues = [ordered]@{
    'B10' = 'WazirX'
    'C10' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'E10' = '9WazirXWRX'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'E11' = '10MandalaExchangeTokenMDX'
    'B12' = 'LiechtensteinCryptoassetsExchange'
    'C12' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'E12' = '11LiechtensteinCryptoassetsExchangeLCX'
    'B13' = 'BitrueCoin'
    'C13' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'E13' = '12BitrueCoinBTR'
    'B14' = 'BitMartToken'
    'C14' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'E14' = '13BitMartTokenBMX'
    'B15' = 'MCDex'
    'C15' = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
    'E15' = '14MCDexMCB'
    'B16' = 'BitForexToken'
    'C16' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'E16' = '15BitForexTokenBF'
    'B17' = 'CoinExToken'
    'C17' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'E17' = '16CoinExTokenCET'
    'B18' = 'One'
    'C18' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'E18' = '17OneONEWorstin24h'
    'E20' = '19BitKanKANBestin24h'
    'E22' = '21NitroExNTX'
    'B42' = 'BKEXToken'
    'C42' = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
    'E42' = '41BKEXTokenBKK'
    'B43' = 'CEJI'
    'C43' = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
    'E43' = '42CEJICEJI'
}

# Price cells hold numeric-looking text (e.g. "243.19") that must stay
# stored as text, exactly like the original inline strings. Force the
# number format to Text before assigning, then restore the cell style so
# no extra formatting is introduced.
$numericTextValues = [ordered]@{
    'D2' = '243.19'
    'D3' = '22.97'
    'D4' = '5.407'
    'D5' = '0.05970'
    'D6' = '3.421'
    'D7' = '6.510'
    'D8' = '0.8110'
    'D9' = '0.9239'
    'D10' = '0.1439'
    'D11' = '0.07383'
    'D12' = '0.03331'
    'D13' = '0.03088'
    'D14' = '0.09356'
    'D15' = '3.848'
    'D16' = '0.001564'
    'D17' = '0.04693'
    'D18' = '0.0005950'
    'D19' = '0.005874'
    'D20' = '0.001262'
    'D21' = '0.004851'
    'D22' = '0.00008000'
    'D23' = '3.570'
    'D25' = '0.3237'
    'D27' = '0.0002340'
    'D40' = '0.03955'
    'D41' = '0.006370'
    'D42' = '0.1074'
    'D43' = '0.002660'
    'D44' = '0.008896'
    'D45' = '0.00005187'
    'D46' = '0.00000000750'
    'D47' = '0.7001'
    'D48' = '0.002141'
    'D49' = '0.00002100'
    'D50' = '0.0002000'
}

foreach ($addr in $textValues.Keys) {
    $ws.Range($addr).Value = $textValues[$addr]
}

foreach ($addr in $numericTextValues.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $numericTextValues[$addr]
    $rng.Style = "Normal"
}
